$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E4 (end time of the 3rd record on 2014-02-17)
$ws.Range("E4").Value = 0.84722222222222221

# Insert new row 5 with a new day entry (2014-02-18), shifting the
# previous "blank separator row" and the summary rows down by one.
$ws.Rows("5:5").Insert()

$ws.Range("A5").Value = 2014
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 18
$ws.Range("D5").Value = 0.33680555555555558
$ws.Range("E5").Value = 0.5
$ws.Range("F5").Formula = "=(E5-D5)*24*60"

# Copy styles from row 4 cells to the corresponding new row 5 cells
$ws.Range("D4").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("E4").Copy()
$ws.Range("E5").PasteSpecial(-4122)
$ws.Range("F4").Copy()
$ws.Range("F5").PasteSpecial(-4122)
$ws.Range("G4").Copy()
$ws.Range("G5").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Fix up the formulas in the summary rows that now live one row lower
$ws.Range("F7").Formula = "=SUM(F2:F6)"
$ws.Range("F8").Formula = "=F7/60"
$ws.Range("F9").Formula = "=F8/38.5"

$ws.Range("G5").Select()
